# Rename the column headers so that they carry the respective input
# file-format version as a suffix instead of the generic "_old" / "_new"
# markers, then wrap the data range in an Excel Table (ListObject) and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10) describe the "old" / FV2404 side of the diff.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2404"
}

# Column K (11) stays "diff".

# Columns L-U (12-21) describe the "new" / FV2410 side of the diff.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2410"
}

# Turn the whole sheet range into a proper Excel Table with the
# (now renamed) header row, matching the shape of the original data.
$tableRange = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table1"

# Freeze the header row (pane split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
